$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.043.82'
$ws.Range("E2").Value = '  -0.84%  '

$ws.Range("D3").Value = '1.822.02'
$ws.Range("E3").Value = '  -0.18%  '

$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.80'
$ws.Range("E5").Value = '  -1.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4495'
$ws.Range("E7").Value = '  +5.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3699'
$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07305'
$ws.Range("E9").Value = '  +0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8576'
$ws.Range("E10").Value = '  -0.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.78'
$ws.Range("E11").Value = '  -1.31%  '

$ws.Range("D12").Value = '1.821.92'
$ws.Range("E12").Value = '  -0.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.645'
$ws.Range("E13").Value = '  -1.08%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.345'
$ws.Range("E14").Value = '  +0.54%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.35'
$ws.Range("E15").Value = '  +4.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07101'
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008792'
$ws.Range("E18").Value = '  -0.94%  '

$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.99'
$ws.Range("E20").Value = '  -0.68%  '

$ws.Range("D21").Value = '27.065.83'
$ws.Range("E21").Value = '  -0.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.172'
$ws.Range("E22").Value = '  +0.57%  '

$ws.Range("E23").Value = '  +0.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.71'
$ws.Range("E25").Value = '  -0.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.235'
$ws.Range("E26").Value = '  +4.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.44'
$ws.Range("E27").Value = '  +0.32%  '

$ws.Range("E28").Value = '  +0.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.68'
$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08873'
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.7566'
$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.185'
$ws.Range("E32").Value = '  -1.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.961'
$ws.Range("E33").Value = '  +4.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.460'
$ws.Range("E34").Value = '  +0.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.000'
$ws.Range("E35").Value = '  -0.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.094'
$ws.Range("E36").Value = '  -1.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01970'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05244'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5330'
$ws.Range("E39").Value = '  +5.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.892'
$ws.Range("E40").Value = '  +0.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.152'
$ws.Range("E41").Value = '  +0.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1706'
$ws.Range("E42").Value = '  +0.54%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5252'
$ws.Range("E43").Value = '  +10.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.526'
$ws.Range("E44").Value = '  -1.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.67'
$ws.Range("E45").Value = '  +0.70%  '

$ws.Range("E46").Value = '  +8.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.62'
$ws.Range("E47").Value = '  -1.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  -0.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.669'
$ws.Range("E49").Value = '  +0.41%  '

$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9211'
$ws.Range("E51").Value = '  +0.46%  '
